# "setting up vegetarian sub"
# Insert a new "broth" sheet (a vegetarian substitute for meat stock) right
# after the "red meat" tab, fill in its property/value rows, and tidy up a
# couple of leftover selection artifacts on the "meat" sheet.

$wb = $excel.ActiveWorkbook

# 1. Add the new sheet immediately after "red meat" (this lands it right
#    before "pork", matching the workbook.xml <sheets> order in the diff).
$redMeat = $wb.Worksheets.Item("red meat")
$broth = $wb.Worksheets.Add($null, $redMeat)
$broth.Name = "broth"

# 2. Populate the property/value table for the new sheet.
$broth.Range("A1").Value = "property"
$broth.Range("B1").Value = "value"

$broth.Range("A2").Value = "name"
$broth.Range("B2").Value = "broth"

$broth.Range("A3").Value = "healthy"
$broth.Range("B3").Value = $true

$broth.Range("A4").Value = "food super group"
$broth.Range("B4").Value = "meat"

$broth.Range("A5").Value = "vegetarian substitute"
$broth.Range("B5").Value = "vegetable broth"

# Match the final selection left on the new sheet.
$null = $broth.Range("B6").Select()

# 3. Tidy up the "meat" sheet's selection (now a full-table selection,
#    A1:B10, instead of the old single-cell A10 selection).
$meat = $wb.Worksheets.Item("meat")
$meat.Activate()
$null = $meat.Range("A1:B10").Select()

# 4. Leave "broth" as the active sheet/tab, matching the new activeTab.
$broth.Activate()
$null = $broth.Range("B6").Select()
